# Results: Add accessibility panel in results section
# - Updates the CO2 label (removes rich "subscript 2" formatting, now plain "CO2")
# - Appends new accessibility-panel related localisation rows to the Labels sheet
# - Moves the active sheet/selection from Choices!E32 to Labels!D20

$wb = $excel.ActiveWorkbook

# --- 1. Fix the CO2 label cells on the Labels sheet (plain text instead of rich "CO" + subscript "2") ---
$labels = $wb.Worksheets.Item("Labels")
$labels.Range("C25").Value = "CO2"
$labels.Range("D25").Value = "CO2"

# --- 2. Append the new accessibility-panel localisation rows (34-45) ---
$newRows = @(
  @{Row=34; B="accessibilityPanel.title";               C="Accessibilité";       D="Accessibility"},
  @{Row=35; B="accessibilityPanel.locationsTitle";       C="Logement";            D="Locations"},
  @{Row=36; B="accessibilityPanel.bothAddresses";        C="Les deux";            D="Both"},
  @{Row=37; B="accessibilityPanel.firstAddressOnly";     C="Logement #1";         D="House #1"},
  @{Row=38; B="accessibilityPanel.secondAddressOnly";    C="Logement #2";         D="House #2"},
  @{Row=39; B="accessibilityPanel.travelTimeTitle";      C="Temps de trajet";     D="Travel time"},
  @{Row=40; B="accessibilityPanel.15min";                C="15 min.";             D="15 min."},
  @{Row=41; B="accessibilityPanel.30min";                C="30 min.";             D="30 min."},
  @{Row=42; B="accessibilityPanel.45min";                C="45 min.";             D="45 min."},
  @{Row=43; B="accessibilityPanel.modeOfTransportTitle"; C="Mode de transport";   D="Mode of transport"},
  @{Row=44; B="accessibilityPanel.minimize";              C="Réduire la fenêtre"; D="Minimize the panel"},
  @{Row=45; B="accessibilityPanel.expand";                C="Agrandir la fenêtre";D="Maximize the panel"}
)

foreach ($item in $newRows) {
  $labels.Cells.Item($item.Row, 1).Value = "results"
  $labels.Cells.Item($item.Row, 2).Value = $item.B
  $labels.Cells.Item($item.Row, 3).Value = $item.C
  $labels.Cells.Item($item.Row, 4).Value = $item.D
}

# --- 3. Move the selection on the Choices sheet back to its existing cell (tab no longer active) ---
$choices = $wb.Worksheets.Item("Choices")
$choices.Range("E32").Select()

# --- 4. Activate the Labels sheet and select D20 (becomes the active tab/selection) ---
$labels.Activate()
$labels.Range("D20").Select()
